$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: cells in column D whose new value is a "pure number" string need to be
# forced to Text so Excel doesn't silently convert them to a numeric type
# (matches the source data, which stores prices as literal text).
$textCells = @("D5","D6","D7","D8","D9","D10","D12","D15","D20","D21","D22","D23","D24","D25","D26","D27","D28","D29","D30","D31","D32","D33","D34","D35","D36","D37","D38","D40","D41","D42","D43","D44","D47","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Row 2 - Bitcoin
$ws.Range("D2").Value = "25.915.83"
$ws.Range("E2").Value = "  -0.34%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.633.58"
$ws.Range("E3").Value = "  -0.51%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  -0.18%  "

# Row 5 - BNB
$ws.Range("D5").Value = "215.93"
$ws.Range("E5").Value = "  +0.49%  "

# Row 6 - XRP
$ws.Range("D6").Value = "0.5117"
$ws.Range("E6").Value = "  +0.47%  "

# Row 7 - USDC
$ws.Range("D7").Value = "1.002"
$ws.Range("E7").Value = "  -0.13%  "

# Row 8 - Cardano
$ws.Range("D8").Value = "0.2567"
$ws.Range("E8").Value = "  +0.04%  "

# Row 9 - Dogecoin
$ws.Range("D9").Value = "0.06348"
$ws.Range("E9").Value = "  -0.26%  "

# Row 10 - Solana
$ws.Range("D10").Value = "19.49"
$ws.Range("E10").Value = "  -0.38%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.19%  "

# Row 12 - Polkadot
$ws.Range("D12").Value = "4.253"
$ws.Range("E12").Value = "  -0.68%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.633.99"
$ws.Range("E13").Value = "  -1.08%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "1.858.59"
$ws.Range("E14").Value = "  -0.55%  "

# Row 15 - Polygon
$ws.Range("D15").Value = "0.5522"
$ws.Range("E15").Value = "  +1.45%  "

# Row 16 - Litecoin
$ws.Range("E16").Value = "  -0.59%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "0.0₅7619"
$ws.Range("E17").Value = "  -1.50%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "25.940.20"
$ws.Range("E18").Value = "  -0.38%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  -0.03%  "

# Row 20 - BitcoinCash
$ws.Range("D20").Value = "195.24"
$ws.Range("E20").Value = "  -0.85%  "

# Row 21 - Uniswap
$ws.Range("D21").Value = "4.427"
$ws.Range("E21").Value = "  +0.17%  "

# Row 22 - Avalanche
$ws.Range("D22").Value = "9.871"
$ws.Range("E22").Value = "  -0.62%  "

# Row 23 - Chainlink
$ws.Range("D23").Value = "6.033"
$ws.Range("E23").Value = "  +0.03%  "

# Row 24 - BinanceUSD
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  -0.29%  "

# Row 25 - Toncoin
$ws.Range("D25").Value = "1.891"
$ws.Range("E25").Value = "  +1.24%  "

# Row 26 - Monero
$ws.Range("D26").Value = "142.34"
$ws.Range("E26").Value = "  +0.49%  "

# Row 27 - Stellar
$ws.Range("D27").Value = "0.1260"
$ws.Range("E27").Value = "  +5.93%  "

# Row 28 - was EthereumClassic, now Cosmos
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "6.769"
$ws.Range("E28").Value = "  -0.81%  "

# Row 29 - was Cosmos, now EthereumClassic
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "15.57"
$ws.Range("E29").Value = "  -0.13%  "

# Row 30 - PancakeSwap
$ws.Range("D30").Value = "1.243"
$ws.Range("E30").Value = "  +0.52%  "

# Row 31 - Hedera
$ws.Range("D31").Value = "0.04919"
$ws.Range("E31").Value = "  +1.22%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("D32").Value = "3.239"
$ws.Range("E32").Value = "  -0.49%  "

# Row 33 - Filecoin
$ws.Range("D33").Value = "3.191"
$ws.Range("E33").Value = "  +0.76%  "

# Row 34 - LidoDAOToken
$ws.Range("D34").Value = "1.548"
$ws.Range("E34").Value = "  +1.52%  "

# Row 35 - HuobiToken
$ws.Range("D35").Value = "2.374"
$ws.Range("E35").Value = "  +0.29%  "

# Row 36 - ARBITRUM
$ws.Range("D36").Value = "0.8984"
$ws.Range("E36").Value = "  +0.17%  "

# Row 37 - ImmutableX
$ws.Range("D37").Value = "0.5533"
$ws.Range("E37").Value = "  +1.46%  "

# Row 38 - MXToken
$ws.Range("D38").Value = "2.539"
$ws.Range("E38").Value = "  -1.58%  "

# Row 39 - Maker
$ws.Range("D39").Value = "1.116.34"
$ws.Range("E39").Value = "  -2.16%  "

# Row 40 - VeChain
$ws.Range("D40").Value = "0.01557"
$ws.Range("E40").Value = "  -0.37%  "

# Row 41 - PaxDollar
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  -0.19%  "

# Row 42 - FraxShare
$ws.Range("D42").Value = "5.589"
$ws.Range("E42").Value = "  +3.47%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").Value = "0.7958"
$ws.Range("E43").Value = "  -1.79%  "

# Row 44 - Quant
$ws.Range("D44").Value = "97.86"
$ws.Range("E44").Value = "  -1.43%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.769.04"
$ws.Range("E45").Value = "  -0.60%  "

# Row 46 - BabyDogeCoin
$ws.Range("E46").Value = "  -7.81%  "

# Row 47 - Mantle
$ws.Range("D47").Value = "0.4436"
$ws.Range("E47").Value = "  -2.14%  "

# Row 48 - Frax
$ws.Range("D48").Value = "1.004"
$ws.Range("E48").Value = "  +0.40%  "

# Row 49 - Aave
$ws.Range("D49").Value = "54.91"
$ws.Range("E49").Value = "  +0.03%  "

# Row 50 - Cronos
$ws.Range("D50").Value = "0.05135"
$ws.Range("E50").Value = "  +1.51%  "

# Row 51 - EnergySwap
$ws.Range("D51").Value = "7.551"
$ws.Range("E51").Value = "  +2.74%  "
